# Insert a new data row at row 87 (pushing the existing rows 87-161 down to
# 88-162), then populate the new row with its values. This mirrors the
# author's edit: a new weekly price observation (2021-09-08 / serial 44447)
# was inserted into the "Acelga" price history kept in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Insert()

$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("D87").Value = 44447
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112009
$ws.Range("G87").Value = "Acelga"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 500
$ws.Range("K87").Value = 2800
$ws.Range("L87").Value = 2800
$ws.Range("M87").Value = 2800
$ws.Range("N87").Value = '$/docena de atados (4 kilos)'
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 700
$ws.Range("Q87").Value = 4
$ws.Range("R87").Value = "Hortaliza"
